# MOM06242015.docx — add "Date:" line, move the _GoBack bookmark, and
# fix up the handful of <w:lastRenderedPageBreak/> markers that drifted
# when Word re-paginated after the new first paragraph was added.

$d = $word.ActiveDocument

function Insert-RawXml($Range, $BodyXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        $BodyXml +
        '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $Range.InsertXML($pkg)
}

function Find-WholeWord($Text) {
    $rng = $d.Content
    $rng.Start = 0
    $rng.End = $d.Content.End
    $ok = $rng.Find.Execute($Text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find-WholeWord: '$Text' not found"
    }
    return $rng
}

# ---------------------------------------------------------------------
# 1) "15" task row — lastRenderedPageBreak now lands here.
# ---------------------------------------------------------------------
$r15 = Find-WholeWord "15"
$body15 = '<w:body><w:p w:rsidR="00A2069C" w:rsidRDefault="00A2069C">' +
    '<w:pPr><w:jc w:val="right"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t>15</w:t></w:r>' +
    '</w:p></w:body>'
Insert-RawXml $r15 $body15

# ---------------------------------------------------------------------
# 2) "aws" row — lastRenderedPageBreak no longer lands here.
# ---------------------------------------------------------------------
$rAws = Find-WholeWord "aws"
$bodyAws = '<w:body><w:p w:rsidR="00A2069C" w:rsidRDefault="00A2069C">' +
    '<w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">Pradeep has details, need to configure based on </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>aws</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p></w:body>'
Insert-RawXml $rAws $bodyAws

# ---------------------------------------------------------------------
# 3) "17" task row — lastRenderedPageBreak no longer lands here.
# ---------------------------------------------------------------------
$r17 = Find-WholeWord "17"
$body17 = '<w:body><w:p w:rsidR="00A2069C" w:rsidRDefault="00A2069C">' +
    '<w:pPr><w:jc w:val="right"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>17</w:t></w:r>' +
    '</w:p></w:body>'
Insert-RawXml $r17 $body17

# ---------------------------------------------------------------------
# 4) "34" task row — lastRenderedPageBreak now lands here.
# ---------------------------------------------------------------------
$r34 = Find-WholeWord "34"
$body34 = '<w:body><w:p w:rsidR="00A2069C" w:rsidRDefault="00A2069C">' +
    '<w:pPr><w:jc w:val="right"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t>34</w:t></w:r>' +
    '</w:p></w:body>'
Insert-RawXml $r34 $body34

# ---------------------------------------------------------------------
# 5) "About document" — drop the _GoBack bookmark that used to sit here.
# ---------------------------------------------------------------------
$rAbout = Find-WholeWord "About document"
$bodyAbout = '<w:body><w:p w:rsidR="00A2069C" w:rsidRDefault="00A2069C">' +
    '<w:pPr><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>About document</w:t></w:r>' +
    '</w:p></w:body>'
Insert-RawXml $rAbout $bodyAbout

# ---------------------------------------------------------------------
# 6) _GoBack bookmark now sits on the blank paragraph right after the
#    first (names) table, instead of next to "About document".
# ---------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$rGoBack = $d.Range($t1.Range.End, $t1.Range.End)
$bodyGoBack = '<w:body><w:p w:rsidR="00B831AB" w:rsidRDefault="00A2069C">' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p></w:body>'
Insert-RawXml $rGoBack $bodyGoBack

# ---------------------------------------------------------------------
# 7) New "Date: 6/24/2015 at 9.30 AM EST" paragraph + blank line, at the
#    very top of the document (before the first table). Inserted last so
#    none of the offsets used above shift underneath us.
# ---------------------------------------------------------------------
$rTop = $d.Range(0, 0)
$bodyTop = '<w:body>' +
    '<w:p><w:pPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>Date: 6/</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>24</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr><w:t>/2015 at 9.30 AM EST</w:t></w:r>' +
    '</w:p>' +
    '<w:p/>' +
    '</w:body>'
Insert-RawXml $rTop $bodyTop

Write-Output "done"
